$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 257, shifting the existing rows 257-265
# down to 258-266 (mirrors the weekly data-entry pattern: a fresh week's
# record is prepended above the rolling history for this market/product).
$ws.Rows.Item(257).Insert()

# Populate the newly inserted row 257 with this week's observation.
$ws.Cells.Item(257, 1).Value = 4
$ws.Cells.Item(257, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(257, 3).Value = "Los Lagos"
$ws.Cells.Item(257, 4).Value = 44747
$ws.Cells.Item(257, 5).Value = 10
$ws.Cells.Item(257, 6).Value = 100112044
$ws.Cells.Item(257, 7).Value = "Perejil"
$ws.Cells.Item(257, 8).Value = "Sin especificar"
$ws.Cells.Item(257, 9).Value = "Primera"
$ws.Cells.Item(257, 10).Value = 180
$ws.Cells.Item(257, 11).Value = 5500
$ws.Cells.Item(257, 12).Value = 6000
$ws.Cells.Item(257, 13).Value = 5750
$ws.Cells.Item(257, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(257, 15).Value = "Región Metropolitana"
$ws.Cells.Item(257, 16).Value = 1917
$ws.Cells.Item(257, 17).Value = 3
$ws.Cells.Item(257, 18).Value = "Hortaliza"
